$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''97.680.30'
$ws.Range("E2").Value = '  -1.15%  '

# Row 3
$ws.Range("D3").Value = '''3.396.76'
$ws.Range("E3").Value = '  +2.43%  '

# Row 4
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").Value = '''253.93'
$ws.Range("E5").Value = '  -0.41%  '

# Row 6
$ws.Range("D6").Value = '''650.49'
$ws.Range("E6").Value = '  +3.36%  '

# Row 7
$ws.Range("E7").Value = '  +1.85%  '

# Row 8
$ws.Range("D8").Value = '''0.432'
$ws.Range("E8").Value = '  +5.70%  '

# Row 9
$ws.Range("D9").Value = '''1.07'
$ws.Range("E9").Value = '  +6.86%  '

# Row 10
$ws.Range("D10").Value = '''1.00'
$ws.Range("E10").Value = '  +0.08%  '

# Row 11
$ws.Range("D11").Value = '''3.401.35'
$ws.Range("E11").Value = '  +2.62%  '

# Row 12
$ws.Range("E12").Value = '  +3.67%  '

# Row 13
$ws.Range("D13").Value = '''41.43'
$ws.Range("E13").Value = '  -4.20%  '

# Row 14
$ws.Range("D14").Value = '''6.39'
$ws.Range("E14").Value = '  +18.31%  '

# Row 15
$ws.Range("D15").Value = '''0.0000259'
$ws.Range("E15").Value = '  +3.10%  '

# Row 16
$ws.Range("D16").Value = '''97.419.54'
$ws.Range("E16").Value = '  -1.11%  '

# Row 17
$ws.Range("D17").Value = '''4.038.61'
$ws.Range("E17").Value = '  +2.62%  '

# Row 18
$ws.Range("D18").Value = '''8.53'
$ws.Range("E18").Value = '  +29.18%  '

# Row 19
$ws.Range("D19").Value = '''3.397.70'
$ws.Range("E19").Value = '  +2.46%  '

# Row 20
$ws.Range("D20").Value = '''17.49'
$ws.Range("E20").Value = '  +8.18%  '

# Row 21
$ws.Range("D21").Value = '''0.498'
$ws.Range("E21").Value = '  +46.25%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''10.75'
$ws.Range("E22").Value = '  +12.47%  '

# Row 23
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").Value = '''3.43'
$ws.Range("E23").Value = '  -1.67%  '

# Row 24
$ws.Range("D24").Value = '''506.48'
$ws.Range("E24").Value = '  +4.19%  '

# Row 25
$ws.Range("E25").Value = '  +0.92%  '

# Row 26
$ws.Range("D26").Value = '''6.18'
$ws.Range("E26").Value = '  +2.66%  '

# Row 27
$ws.Range("D27").Value = '''98.92'
$ws.Range("E27").Value = '  +8.77%  '

# Row 28
$ws.Range("D28").Value = '''12.73'
$ws.Range("E28").Value = '  +3.49%  '

# Row 29
$ws.Range("D29").Value = '''3.587.25'
$ws.Range("E29").Value = '  +2.71%  '

# Row 30
$ws.Range("D30").Value = '''0.152'
$ws.Range("E30").Value = '  +2.96%  '

# Row 31
$ws.Range("D31").Value = '''0.205'
$ws.Range("E31").Value = '  +6.83%  '

# Row 32
$ws.Range("D32").Value = '''11.37'
$ws.Range("E32").Value = '  +3.53%  '

# Row 33
$ws.Range("D33").Value = '''0.998'
$ws.Range("E33").Value = '  -0.19%  '

# Row 34
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  -0.04%  '

# Row 35
$ws.Range("D35").Value = '''0.567'
$ws.Range("E35").Value = '  +16.60%  '

# Row 36
$ws.Range("D36").Value = '''29.69'
$ws.Range("E36").Value = '  +5.22%  '

# Row 37
$ws.Range("D37").Value = '''2.26'
$ws.Range("E37").Value = '  +14.21%  '

# Row 38
$ws.Range("D38").Value = '''7.68'
$ws.Range("E38").Value = '  +3.35%  '

# Row 39
$ws.Range("D39").Value = '''526.05'
$ws.Range("E39").Value = '  +4.84%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.153'
$ws.Range("E40").Value = '  +1.09%  '

# Row 41
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '''1.42'
$ws.Range("E41").Value = '  +12.25%  '

# Row 42
$ws.Range("D42").Value = '''24.73'
$ws.Range("E42").Value = '  -0.05%  '

# Row 43
$ws.Range("D43").Value = '''0.857'
$ws.Range("E43").Value = '  +7.04%  '

# Row 44
$ws.Range("E44").Value = '  -4.55%  '

# Row 45
$ws.Range("D45").Value = '''0.0419'
$ws.Range("E45").Value = '  +18.74%  '

# Row 46
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '''5.47'
$ws.Range("E46").Value = '  +13.36%  '

# Row 47
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''3.28'
$ws.Range("E47").Value = '  +3.73%  '

# Row 48
$ws.Range("D48").Value = '''8.23'
$ws.Range("E48").Value = '  +10.62%  '

# Row 49
$ws.Range("E49").Value = '  -0.02%  '

# Row 50
$ws.Range("E50").Value = '  +10.69%  '

# Row 51
$ws.Range("D51").Value = '''2.05'
$ws.Range("E51").Value = '  +3.32%  '
